$d = $word.ActiveDocument

$d.Content.Find.Execute("181÷3=60, 1", $true, $false, $false, $false, $false, $true, 1, $false, "869÷4=217, 1", 2)
$d.Content.Find.Execute("716÷5=143, 1", $true, $false, $false, $false, $false, $true, 1, $false, "806÷4=201, 2", 2)
$d.Content.Find.Execute("245÷6=40, 5", $true, $false, $false, $false, $false, $true, 1, $false, "835÷7=119, 2", 2)
$d.Content.Find.Execute("347÷7=49, 4", $true, $false, $false, $false, $false, $true, 1, $false, "147÷4=36, 3", 2)
$d.Content.Find.Execute("670÷2=335, 0", $true, $false, $false, $false, $false, $true, 1, $false, "937÷8=117, 1", 2)
$d.Content.Find.Execute("607÷2=303, 1", $true, $false, $false, $false, $false, $true, 1, $false, "126÷4=31, 2", 2)
$d.Content.Find.Execute("635÷7=90, 5", $true, $false, $false, $false, $false, $true, 1, $false, "119÷5=23, 4", 2)
$d.Content.Find.Execute("547÷7=78, 1", $true, $false, $false, $false, $false, $true, 1, $false, "104÷5=20, 4", 2)
$d.Content.Find.Execute("499÷7=71, 2", $true, $false, $false, $false, $false, $true, 1, $false, "968÷7=138, 2", 2)
$d.Content.Find.Execute("669÷8=83, 5", $true, $false, $false, $false, $false, $true, 1, $false, "956÷9=106, 2", 2)
$d.Content.Find.Execute("820÷8=102, 4", $true, $false, $false, $false, $false, $true, 1, $false, "192÷5=38, 2", 2)
$d.Content.Find.Execute("826÷9=91, 7", $true, $false, $false, $false, $false, $true, 1, $false, "171÷9=19, 0", 2)
$d.Content.Find.Execute("426÷3=142, 0", $true, $false, $false, $false, $false, $true, 1, $false, "522÷3=174, 0", 2)
$d.Content.Find.Execute("940÷3=313, 1", $true, $false, $false, $false, $false, $true, 1, $false, "346÷4=86, 2", 2)
$d.Content.Find.Execute("502÷8=62, 6", $true, $false, $false, $false, $false, $true, 1, $false, "685÷2=342, 1", 2)
$d.Content.Find.Execute("586÷9=65, 1", $true, $false, $false, $false, $false, $true, 1, $false, "850÷2=425, 0", 2)
$d.Content.Find.Execute("565÷3=188, 1", $true, $false, $false, $false, $false, $true, 1, $false, "879÷4=219, 3", 2)
$d.Content.Find.Execute("546÷3=182, 0", $true, $false, $false, $false, $false, $true, 1, $false, "522÷7=74, 4", 2)
$d.Content.Find.Execute("763÷2=381, 1", $true, $false, $false, $false, $false, $true, 1, $false, "985÷6=164, 1", 2)
$d.Content.Find.Execute("314÷2=157, 0", $true, $false, $false, $false, $false, $true, 1, $false, "621÷2=310, 1", 2)
$d.Content.Find.Execute("871÷2=435, 1", $true, $false, $false, $false, $false, $true, 1, $false, "179÷9=19, 8", 2)
$d.Content.Find.Execute("538÷5=107, 3", $true, $false, $false, $false, $false, $true, 1, $false, "261÷6=43, 3", 2)
$d.Content.Find.Execute("733÷4=183, 1", $true, $false, $false, $false, $false, $true, 1, $false, "656÷6=109, 2", 2)
$d.Content.Find.Execute("259÷4=64, 3", $true, $false, $false, $false, $false, $true, 1, $false, "713÷2=356, 1", 2)
$d.Content.Find.Execute("301÷2=150, 1", $true, $false, $false, $false, $false, $true, 1, $false, "932÷9=103, 5", 2)
